$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item("Прямоугольник 6")
$tr = $sh.TextFrame.TextRange

# Replace the whole original run's text with the new, longer wording.
$full = $tr.Characters(1, 193)
$full.Text = "Для меня эта тема актуальна так как газы очень часто встречаются и проблема экологической обстановки в мире на сегодняшний день очень плачевна я хочу на примере показать что происходит с воздухом при выбросе разных химических, ядовитых газов.  "

# Force the text to split into five runs (matching the authored edit) by
# touching each sub-range with a formatting no-op (it is already italic).
$r1 = $tr.Characters(1, 32)
$r1.Font.Italic = $true
$r2 = $tr.Characters(33, 4)
$r2.Font.Italic = $true
$r3 = $tr.Characters(37, 4)
$r3.Font.Italic = $true
$r4 = $tr.Characters(41, 175)
$r4.Font.Italic = $true
$r5 = $tr.Characters(216, 29)
$r5.Font.Italic = $true

# The textbox auto-fits to its text; grow it to the new rendered height.
$sh.Height = 312.62346456692916
